# Weekly update: two new Chirimoya price records were entered at the top of
# the "Macroferia Regional de Talca" data block (rows 141-142), pushing all
# the existing records for that block down by two rows (old row 141 -> 143,
# ..., old row 186 -> 188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 141; everything that used
# to live in rows 141..186 shifts down to rows 143..188.
$ws.Rows("141:142").Insert()

# Populate the two new rows with the new weekly entries. Columns A,B,C,E,F,
# G,H,I,J,K,Q,R,T are constant for every record in this subset.
$newRows = @(
    @{ Row = 141; D = 45229; L = "Primera"; M = 280; N = 20000; O = 20000; P = 20000; S = 2000 },
    @{ Row = 142; D = 45229; L = "Segunda"; M = 200; N = 18000; O = 18000; P = 18000; S = 1800 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 10
}
